# Fixing input data for some tests.
$wb = $excel.ActiveWorkbook

$wsScenarios   = $wb.Worksheets.Item("_set_scenarios")
$wsUncertainty = $wb.Worksheets.Item("_set_uncertainty")
$wsProducts    = $wb.Worksheets.Item("_set_products")

# Rename header cells on each of the "set" sheets.
$wsProducts.Range("A1").Value = "p_Names"
$wsUncertainty.Range("A1").Value = "u_Names"
$wsScenarios.Range("A1").Value = "s_Names"

# Make "_set_scenarios" the active sheet/tab (was "_set_products").
$wsScenarios.Activate()
$wsScenarios.Range("D15").Select()

$wsUncertainty.Range("A2").Select()

$wsProducts.Range("F37").Select()

$wsScenarios.Activate()
